# #5: property aircraft done
#
# The "property_category" column on each per-category sheet ("土地"/land,
# "建物"/building, "汽車"/car, ...) had been stamped with the wrong
# category label (several sheets were incorrectly showing "land").
# Fix the per-sheet "property_category" values so each sheet reports its
# own category:
#   - 建物 (building) sheet: rows 2-5, column I -> "building"
#   - 汽車 (car) sheet:      row 2,   column H -> "car"

$wb = $excel.ActiveWorkbook

# 建物 (building) sheet - 3rd column group ("property_category") is column I
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"
$wsBuilding.Range("I4").Value = "building"
$wsBuilding.Range("I5").Value = "building"

# 汽車 (car) sheet - "property_category" is column H
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
